# Fix Training Data Issue
# The "Date" column (BF) was populated with a malformed value
# ("4-18-2012-13") combining month-day with the two-year season label.
# Replace it with the correct ISO date "2013-04-18" for every data row.
#
# Assigning a date-looking string straight to .Value2 makes Excel parse it
# as a serial date, so we stage the literal text in a scratch cell that is
# explicitly formatted as Text ("@"), copy it, and paste-special only the
# values into the target cells -- that keeps the destination cells holding
# the literal string "2013-04-18" (not a date serial) without leaving any
# stray number-format applied to them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "4-18-2012-13"
$newValue = "2013-04-18"

# Scratch cell inside the sheet's existing used range (A1 is blank) so we
# don't grow the worksheet dimensions.
$scratch = $ws.Range("A1")

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)  # column BF
    if ($cell.Value2 -eq $oldValue) {
        $scratch.NumberFormat = "@"
        $scratch.Value2 = $newValue
        $scratch.Copy()
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}

$scratch.Clear()
$excel.CutCopyMode = $false
